$d = $word.ActiveDocument

$pairs = @(
    @("886÷3=295, 1", "117÷2=58, 1"),
    @("669÷6=111, 3", "268÷4=67, 0"),
    @("635÷8=79, 3", "785÷3=261, 2"),
    @("978÷8=122, 2", "224÷2=112, 0"),
    @("123÷9=13, 6", "255÷7=36, 3"),
    @("121÷3=40, 1", "709÷8=88, 5"),
    @("524÷3=174, 2", "440÷5=88, 0"),
    @("174÷3=58, 0", "293÷3=97, 2"),
    @("586÷6=97, 4", "467÷7=66, 5"),
    @("518÷6=86, 2", "337÷3=112, 1"),
    @("684÷9=76, 0", "314÷3=104, 2"),
    @("691÷5=138, 1", "848÷5=169, 3"),
    @("172÷6=28, 4", "975÷9=108, 3"),
    @("797÷6=132, 5", "429÷7=61, 2"),
    @("301÷6=50, 1", "996÷8=124, 4"),
    @("374÷2=187, 0", "597÷4=149, 1"),
    @("365÷4=91, 1", "470÷4=117, 2"),
    @("506÷5=101, 1", "598÷4=149, 2"),
    @("387÷5=77, 2", "675÷4=168, 3"),
    @("390÷8=48, 6", "798÷7=114, 0"),
    @("610÷7=87, 1", "902÷5=180, 2"),
    @("933÷2=466, 1", "747÷7=106, 5"),
    @("211÷8=26, 3", "190÷4=47, 2"),
    @("569÷4=142, 1", "297÷4=74, 1"),
    @("770÷3=256, 2", "424÷3=141, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
